# Publish IG 1.0.1
$wb = $excel.ActiveWorkbook

# 1. Rename the "Include from unknown" sheet to "Include #0"
$ws2 = $wb.Worksheets.Item("Include from unknown")
$ws2.Name = "Include #0"

# 2. Metadata sheet updates
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value 1.0.0 -> 1.0.1
$ws1.Range("B3").Value = "1.0.1"

# Update Contact value
$ws1.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10), before "Description" (row 11)
$ws1.Rows.Item(11).Insert()

# Copy formatting from the row above so the new row matches the table style
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
